$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 135
$ws.Range("F3").Value = 356
$ws.Range("F5").Value = 13
$ws.Range("F6").Value = 1234
$ws.Range("G6").Value = 156
$ws.Range("F7").Value = 444
$ws.Range("F9").Value = 172
$ws.Range("F10").Value = 148
$ws.Range("F11").Value = 1041
$ws.Range("F13").Value = 261
$ws.Range("F14").Value = 163
$ws.Range("F16").Value = 1468
$ws.Range("F17").Value = 540
$ws.Range("F19").Value = 337
$ws.Range("F21").Value = 795
$ws.Range("F22").Value = 1138
$ws.Range("F23").Value = 57
$ws.Range("F25").Value = 2633
$ws.Range("F26").Value = 1404
$ws.Range("F28").Value = 26
$ws.Range("F29").Value = 359
$ws.Range("F30").Value = 403
$ws.Range("F31").Value = 1137
$ws.Range("F32").Value = 807
$ws.Range("F33").Value = 1279
$ws.Range("F34").Value = 151
$ws.Range("F37").Value = 559
$ws.Range("F38").Value = 649
$ws.Range("F39").Value = 823
$ws.Range("F40").Value = 348

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 178
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 2
$ws.Range("F15").Value = 609
$ws.Range("F16").Value = 22

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 867

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 867
$ws.Range("F5").Value = 135
$ws.Range("F6").Value = 356
$ws.Range("F9").Value = 178
$ws.Range("F10").Value = 1234
$ws.Range("G10").Value = 156
$ws.Range("F11").Value = 444
$ws.Range("F13").Value = 172
$ws.Range("F15").Value = 148
$ws.Range("F17").Value = 261
$ws.Range("F19").Value = 163
$ws.Range("F21").Value = 1468
$ws.Range("F22").Value = 540
$ws.Range("F24").Value = 337
$ws.Range("F26").Value = 1138
$ws.Range("F27").Value = 2633
$ws.Range("F29").Value = 1404
$ws.Range("F32").Value = 26
$ws.Range("F34").Value = 359
$ws.Range("F35").Value = 403
$ws.Range("F36").Value = 1137
$ws.Range("F39").Value = 807
$ws.Range("F40").Value = 1279
$ws.Range("F42").Value = 560
$ws.Range("F43").Value = 649
$ws.Range("F44").Value = 823
$ws.Range("F45").Value = 348
